# EkA RFID Anten projesinin kodu güncellendi. URL uzantısı eklendi.
# Row 8 (A8=7): product code + model type updated to the new "EKA" variant,
# and a new hyperlink with the updated GitHub URL is added to K8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the product code (Ürün Kodu) and model type (Model Türü) for row 8
$ws.Range("B8").Value = "AC-RFD-EK-135-ANT-H4B0-01"
$ws.Range("D8").Value = "EKA"

# Add the hyperlink for K8, pointing to the new repository URL
$ws.Hyperlinks.Add($ws.Range("K8"), "https://github.com/btk42/AC-RFD-EK-135-ANT-H4B0-01", [Type]::Missing, [Type]::Missing, "https://github.com/btk42/AC-RFD-EK-135-ANT-H4B0-01")

# Match the visual style used by the other hyperlink cells (K3, K5, K7)
$ws.Range("K7").Copy()
$ws.Range("K8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update the active cell selection to reflect the author's final cursor position
$ws.Range("N10").Select()
